# Updated symbol list on Wed Dec 14 10:29:21 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as text (not numbers) in the
# workbook, so numeric-looking replacements are written with a leading
# apostrophe to force Excel to keep them as literal text (preserving
# trailing zeros / exact formatting, e.g. "0.06240", "0.0001500").
# The "Volume(1h)" column (E) cells that changed are plain text already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'274.27"
$ws.Range("D3").Value  = "'22.95"
$ws.Range("D4").Value  = "'6.347"
$ws.Range("D5").Value  = "'0.06240"
$ws.Range("D6").Value  = "'3.657"
$ws.Range("D7").Value  = "'6.710"
$ws.Range("D8").Value  = "'1.371"
$ws.Range("D9").Value  = "'0.8312"
$ws.Range("D10").Value = "'0.01375"
$ws.Range("D11").Value = "'0.1636"
$ws.Range("D12").Value = "'0.08298"
$ws.Range("D13").Value = "'0.03369"
$ws.Range("D14").Value = "'0.03106"
$ws.Range("D15").Value = "'0.09319"
$ws.Range("D16").Value = "'3.870"
$ws.Range("D17").Value = "'0.001641"
$ws.Range("D18").Value = "'0.04792"
$ws.Range("D19").Value = "'0.006347"
$ws.Range("D20").Value = "'0.005563"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("D21").Value = "'0.001091"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.729"
$ws.Range("D27").Value = "'0.0002680"
$ws.Range("D40").Value = "'0.04702"
$ws.Range("D41").Value = "'0.007030"
$ws.Range("D42").Value = "'0.1166"
$ws.Range("D43").Value = "'0.003456"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D44").Value = "'0.01214"
$ws.Range("D45").Value = "'0.00006253"
$ws.Range("D47").Value = "'0.9000"
$ws.Range("D48").Value = "'0.02771"
$ws.Range("D49").Value = "'0.00002300"
$ws.Range("D50").Value = "'0.01240"
